# v1.0.1: update TODO goals - the "Example" column (D) for the default-value
# rows now reads "N/A" instead of the stale "[CodePen]()" placeholder.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4-8 in column D get the new shared text "N/A" (rows 2-3 keep their
# original "[CodePen]()" value untouched).
$ws.Range("D4:D8").Value = "N/A"

# Move/record the active selection as it was left in the authored workbook.
$ws.Range("D12").Select()
